# Update "丽水-漫展信息.xlsx" to match the 456a3b4 gh-pages output:
#  - F2 (想去人数 for the HP国风 event) bumps from 502 to 504
#  - A new event ("丽水·逆光ZERO动漫游戏展", 2024-08-03) is inserted as the
#    new row 5, pushing the previously-existing rows 5-8 down to rows 6-9
#  - The "想去人数" for the R动漫嘉年华 row (now row 8) ticks from 1 to 2
#
# Both the "展览" sheet and the "全部类型" sheet carry identical tables and
# both need the same edit.

$wb = $excel.ActiveWorkbook
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # 1) F2: 502 -> 504
    $ws.Range("F2").Value = 504

    # 2) Insert a brand-new row 5 (shifts old rows 5-8 down to 6-9,
    #    dragging their formatting/values along for the ride).
    $ws.Rows.Item(5).Insert()

    # Copy column-A's numbered-index style (border + bold + centering) down
    # into the freshly inserted, still-blank A5 so it matches the rest of
    # the column instead of picking up Excel's auto-inherited format.
    $ws.Range("A6").Copy()
    $ws.Range("A5").PasteSpecial(-4122)

    # 3) Populate the new row 5 with the 逆光ZERO event.
    $ws.Range("A5").Value = 4

    # The date-like "2024-08-03" text would otherwise be auto-parsed into a
    # real Excel date serial; force text storage (matching every other date
    # cell in this column, which are plain strings) and then drop the
    # resulting "@" number-format back off the cell so it stays styleless.
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("B5").Value = "2024-08-03"
    $ws.Range("B5").ClearFormats()

    $ws.Range("C5").Value = "丽水·逆光ZERO动漫游戏展"
    $ws.Range("D5").Value = "丽阳街651号 丽水华侨君澜大饭店"
    $ws.Range("E5").Value = "2024.08.03 10:00-08.03 17:00"
    $ws.Range("F5").Value = 0
    $ws.Range("G5").Value = 40
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=89687"
    $ws.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202407/NMYT1LRl1721639164353.jpeg"

    # 4) The shifted-down rows (6-9) kept their old literal index numbers
    #    (4,5,6,7) in column A - bump each by one so they read 5,6,7,8.
    $ws.Range("A6").Value = 5
    $ws.Range("A7").Value = 6
    $ws.Range("A8").Value = 7
    $ws.Range("A9").Value = 8

    # 5) R动漫嘉年华 (now row 8): 想去人数 1 -> 2
    $ws.Range("F8").Value = 2
}
